$d = $word.ActiveDocument

# 1. " our school," -> " our school name,"
$d.Content.Find.Execute(" our school,", $true, $false, $false, $false, $false,
                         $true, 1, $false, " our school name,", 2) | Out-Null

# 2. Append new sentence after "Argumented Reality and Virtual Reality."
$d.Content.Find.Execute("Argumented Reality and Virtual Reality.", $true, $false, $false, $false, $false,
                         $true, 1, $false,
                         "Argumented Reality and Virtual Reality. I have also added appropriate meta tag, including keywords and description for the club website.",
                         2) | Out-Null

# 3. Append new sentence after "...reduce CSS size."
$d.Content.Find.Execute("I keep the same style for recent events and resources to increase repetition and reduce CSS size.",
                         $true, $false, $false, $false, $false,
                         $true, 1, $false,
                         "I keep the same style for recent events and resources to increase repetition and reduce CSS size. The layout for individual event and resource including a date, title for the item, description for the item, and an image for the item.",
                         2) | Out-Null

# 4. Replace the "All the images..." sentence with the expanded paragraph text.
$d.Content.Find.Execute("All the images have alt text as Alt Text and Writing Great Alt Text articles mentioned.",
                         $true, $false, $false, $false, $false,
                         $true, 1, $false,
                         "I also use the heading tag cautiously, the title for the page uses h1 tag, recent events and resources uses h2 tag, the individual event and resource uses h3 tag for a clearer structure. I have also uses some hew HTML tags to increase semantics, including section and time. Moreover, all the images have alt text as Alt Text and Writing Great Alt Text articles mentioned.",
                         2) | Out-Null
